# Update gh-pages to output generated at 456a3b4
# Applies numeric updates to the F (and one G) column cells across the
# four worksheets of the workbook, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 197
$ws1.Range("F4").Value  = 1130
$ws1.Range("F8").Value  = 13127
$ws1.Range("F9").Value  = 2267
$ws1.Range("G11").Value = 90
$ws1.Range("F12").Value = 53913
$ws1.Range("F15").Value = 301
$ws1.Range("F17").Value = 702
$ws1.Range("F18").Value = 353
$ws1.Range("F19").Value = 2967
$ws1.Range("F20").Value = 851
$ws1.Range("F21").Value = 5130
$ws1.Range("F22").Value = 1244
$ws1.Range("F27").Value = 372
$ws1.Range("F28").Value = 1189
$ws1.Range("F32").Value = 327
$ws1.Range("F35").Value = 59
$ws1.Range("F36").Value = 42
$ws1.Range("F37").Value = 4691
$ws1.Range("F39").Value = 4727
$ws1.Range("F40").Value = 8690
$ws1.Range("F42").Value = 147
$ws1.Range("F43").Value = 115
$ws1.Range("F44").Value = 201
$ws1.Range("F45").Value = 403
$ws1.Range("F47").Value = 68

# --- Sheet 2: 演出 -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 88
$ws2.Range("F12").Value = 1111

# --- Sheet 3: 本地生活 ----------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 779
$ws3.Range("F3").Value = 538
$ws3.Range("F5").Value = 29

# --- Sheet 4: 全部类型 ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 779
$ws4.Range("F3").Value  = 538
$ws4.Range("F4").Value  = 197
$ws4.Range("F6").Value  = 1130
$ws4.Range("F8").Value  = 13127
$ws4.Range("F9").Value  = 13127
$ws4.Range("F10").Value = 2267
$ws4.Range("F14").Value = 702
$ws4.Range("F15").Value = 353
$ws4.Range("F16").Value = 2967
$ws4.Range("F17").Value = 851
$ws4.Range("F18").Value = 88
$ws4.Range("F19").Value = 1244
$ws4.Range("F20").Value = 29
$ws4.Range("F24").Value = 372
$ws4.Range("F26").Value = 1189
$ws4.Range("F31").Value = 327
$ws4.Range("F33").Value = 42
$ws4.Range("F34").Value = 4691
$ws4.Range("F35").Value = 4727
$ws4.Range("F36").Value = 8690
$ws4.Range("F38").Value = 147
$ws4.Range("F39").Value = 201
